$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A's new cell to be stored as text (not auto-parsed as a date)
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = "2025-09-02"
# Restore default styling so the new row matches the other data rows (no explicit style index)
$ws.Range("A18").Style = "Normal"

$ws.Range("B18").Value = 58.25
$ws.Range("C18").Value = 684.4000244140625
$ws.Range("D18").Value = 322.3999938964844
